# Update "Forecast Comparison" sheet with a new Week_Start_Date column
# and corrected week labels / is_holiday_week typing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the old "ASIN" column (B) - this shifts
# ASIN..is_holiday_week one column to the right (B:I -> C:J).
$ws.Columns.Item(2).Insert()

# New header + data for the inserted "Week_Start_Date" column.
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

$weekStarts = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekStarts.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 2)
    # Force text storage so the date keeps its literal "YYYY-MM-DD" form
    # instead of being auto-converted into a date serial number.
    $cell.NumberFormat = "@"
    $cell.Value = $weekStarts[$i]
}

# Correct the "Week" labels in column A: drop the leading zero (W01 -> W1 ... W09 -> W9).
$weekLabels = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
for ($i = 0; $i -lt $weekLabels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $weekLabels[$i]
}

# is_holiday_week moved from column I to column J, and should now be a
# genuine boolean (t="b") rather than a numeric 0/1.
for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 10).Value = $false
}
